$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 78, shifting existing rows 78-148 down to 79-149.
$ws.Range("A78").EntireRow.Insert()

# Populate the newly inserted row 78 with its data (mirrors the layout of the
# surrounding rows, which are unchanged by this edit aside from the shift).
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(78, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(78, 4).Value = 44907
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = 100112042
$ws.Cells.Item(78, 7).Value = "Locoto"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 120
$ws.Cells.Item(78, 11).Value = 17000
$ws.Cells.Item(78, 12).Value = 18000
$ws.Cells.Item(78, 13).Value = 17500
$ws.Cells.Item(78, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(78, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 16).Value = 875
$ws.Cells.Item(78, 17).Value = 20
$ws.Cells.Item(78, 18).Value = "Hortaliza"
